$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A9").Formula2 = "=UNIQUE(A2:A6)"
